$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: Insert the new "6. Operators" / "7. Closures" sections before the
# existing "Example types:" heading, and remove the two blank paragraphs
# that used to sit directly above that heading.
# ---------------------------------------------------------------------------

$xmlBlock = @'
<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>6. Operators</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Style rules for operators are grouped by arity (the number of operands they take).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>When space is permitted around an operator, multiple spaces MAY be used for readability purposes.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>All operators not described here are left undefined.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>The increment/decrement operators MUST NOT have any space between the operator and operand.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Type casting operators MUST NOT have any space within the parentheses</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>All binary arithmetic, comparison, assignment, bitwise, logical, string, and type operators MUST be preceded and followed by at least one space</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>7. Closures</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Closures MUST be declared with a space after the function keyword, and a space before and after the use keyword.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>The opening brace MUST go on the same line, and the closing brace MUST go on the next line following the body.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>There MUST NOT be a space after the opening parenthesis of the argument list or variable list, and there MUST NOT be a space before the closing parenthesis of the argument list or variable list.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>In the argument list and variable list, there MUST NOT be a space before each comma, and there MUST be one space after each comma.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Closure arguments with default values MUST go at the end of the argument list.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>If a return type is present, it MUST follow the same rules as with normal functions and methods; if the use keyword is present, the colon MUST follow the use list closing parentheses with no spaces between the two characters.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p>
'@

$rng = $d.Content
$found = $rng.Find.Execute("Example types:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Example types:' heading" }

$examplePara = $d.Range($rng.Start, $rng.Start).Paragraphs(1)
$prevTwo = $examplePara.Previous(2)

# Delete the two blank paragraphs one paragraph-mark at a time -- deleting a
# range that spans more than one paragraph mark in a single call only
# collapses one of them, so we do it in two single-mark passes instead.
$markOne = $d.Range($prevTwo.Range.Start, $prevTwo.Range.End)
Write-Host "Removing blank paragraph 1, chars=" $markOne.Text.Length
$markOne.Delete()

$rng = $d.Content
$found = $rng.Find.Execute("Example types:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not re-find 'Example types:' heading (pass 2)" }
$examplePara = $d.Range($rng.Start, $rng.Start).Paragraphs(1)
$prevOne = $examplePara.Previous(1)
$markTwo = $d.Range($prevOne.Range.Start, $prevOne.Range.End)
Write-Host "Removing blank paragraph 2, chars=" $markTwo.Text.Length
$markTwo.Delete()

# Re-locate the heading (its start offset moved after the delete) and insert
# the new block of paragraphs immediately before it.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Example types:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not re-find 'Example types:' heading" }

$insertPoint = $d.Range($rng2.Start, $rng2.Start)
$insertPoint.InsertXML($xmlBlock)

Write-Host "Inserted Operators/Closures sections."

# ---------------------------------------------------------------------------
# Part 2: Add a lastRenderedPageBreak before the "Class constant must be
# declared..." run.
# ---------------------------------------------------------------------------

$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Class constant must be declared in uppercase with underscore", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Could not find 'Class constant must be declared' run" }

$breakPoint = $d.Range($rng3.Start, $rng3.Start)
$breakPoint.InsertXML('<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:lastRenderedPageBreak/></w:r>')

Write-Host "Inserted lastRenderedPageBreak."
